$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage values that are safe to assign directly
# (Excel will not auto-convert these to numbers because they contain
# letters, spaces, or non-ASCII digits).
$ws.Range('E3').Value = '  -4.89%  '
$ws.Range('E5').Value = '  -2.28%  '
$ws.Range('E6').Value = '  +1.16%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -1.83%  '
$ws.Range('E9').Value = '  -4.91%  '
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('E12').Value = '  -4.71%  '
$ws.Range('B13').Value = 'Toncoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E13').Value = '  -2.33%  '
$ws.Range('E14').Value = '  -3.73%  '
$ws.Range('E16').Value = '  -3.65%  '
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('E18').Value = '  -4.76%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E19').Value = '  -3.48%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E20').Value = '  -6.82%  '
$ws.Range('E21').Value = '  -2.99%  '
$ws.Range('E22').Value = '  -3.35%  '
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  -4.56%  '
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('E29').Value = '  -5.11%  '
$ws.Range('D30').Value = '0.0₃0976'
$ws.Range('E30').Value = '  -3.80%  '
$ws.Range('E31').Value = '  -4.94%  '
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('E33').Value = '  -3.52%  '
$ws.Range('E34').Value = '  -3.61%  '
$ws.Range('E35').Value = '  -4.95%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -3.33%  '
$ws.Range('E38').Value = '  -0.44%  '
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('E42').Value = '  -1.70%  '
$ws.Range('E43').Value = '  -3.23%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('E45').Value = '  +1.66%  '
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('E47').Value = '  -3.01%  '
$ws.Range('E48').Value = '  -4.45%  '
$ws.Range('E49').Value = '  -3.67%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  -11.18%  '

# Numeric-looking text values in column D must be forced to text so Excel
# does not convert them to actual numbers (these display like
# "66.518.23" which is not a valid Excel number). We temporarily set the
# cell to a text number format, assign the value, then restore the
# original "Normal" style so no stray formatting is left behind.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.518.23'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.512.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.526'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.511.47'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.139'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.12'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.73'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.957.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.379.66'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.509.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.27'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '348.01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.20'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '524.22'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.131'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.64'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.39'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.09'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.37'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '148.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.558'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.69'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.72'
$ws.Range('D50').Style = 'Normal'
